$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '29.634.75'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -2.57%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.005.99'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -4.94%  '
$ws.Range('E4').Value = '  +0.76%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '331.95'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -3.92%  '
$ws.Range('E6').Value = '  +0.58%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5025'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -3.98%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.4267'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -4.14%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '54.91'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.01%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.09163'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -2.41%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '1.128'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -3.95%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '23.57'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -5.63%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '8.138'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -6.57%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '1.984.67'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -3.41%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '6.547'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -5.89%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '95.42'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -6.55%  '
$ws.Range('E17').Value = '  +0.64%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.00001123'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -3.55%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06674'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -0.81%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '19.92'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -6.19%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '1.012'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.63%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '5.992'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -5.62%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '29.650.72'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -2.60%  '
$ws.Range('E24').Value = '  -4.67%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.283'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -0.78%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '159.22'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -2.30%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '20.82'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -5.59%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '6.420'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -6.13%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.319'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -8.38%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '128.88'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -3.87%  '
$ws.Range('E31').Value = '  -7.53%  '
$ws.Range('E32').Value = '  -9.34%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.09961'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -5.68%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '5.862'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -6.56%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '3.804'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -3.20%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '9.554'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -7.86%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.02483'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -5.49%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '1.320'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -3.27%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.06393'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -5.85%  '
$ws.Range('E40').Value = '  -6.45%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '11.78'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -6.48%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.2074'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -6.83%  '
$ws.Range('E43').Value = '  +0.52%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.6378'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -7.08%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '13.65'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -5.83%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '1.289'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -7.80%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '3.532'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.00000000335'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -3.73%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.07003'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -3.50%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.132'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -6.45%  '
